$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ArtifactSetDataTable")

# Update the SetEffects text for "Pale Flame" (row 3, column C):
# MaxNumStacks was 0, now becomes 2.
$ws.Range("C3").Value = ",`n[StatisticBoost;DamageModifier_PhysicalDamage=0.25;],`n,`n[StatisticBoost,Stackable;Attack_Percentage=0.09,MaxNumStacks=2,InitialNumStacks=0;,StatisticBoost,Triggerable;DamageModifier_PhysicalDamage=0.25;],"

# Move the selection to C13 as in the final saved state.
$ws.Range("C13").Select()
